# Weekly driver report update for 2025-04-29
# Updates the "Driver Summary" sheet:
#  - refresh Bad Drivers table metrics (rows 3-5)
#  - collapse the now-empty "Good Drivers" table (row 12) down to a
#    single "No good drivers found." message, dropping the stale
#    trailing blank rows (13-17)
#  - narrow columns B and E to match the new, simpler layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table: updated Critical Minutes / Good Roaming % ---
$ws.Range("C3").Value = 754
$ws.Range("D3").Value = 95
$ws.Range("C4").Value = 866
$ws.Range("D4").Value = 96
$ws.Range("C5").Value = 1620

# --- Good Drivers table (row 12 onward) collapses to a single note ---
$ws.Range("A12:E12").Clear()
$ws.Range("A12").Value = "No good drivers found."

# Any previously-empty trailing rows (13-17) are dropped automatically
# once row 12 no longer spills past column E - nothing left to delete.
$ws.Range("A13:J17").Clear()

# Re-touch a cell in column J of the still-active row so the sheet's
# used range keeps spanning out to column J (matching the original
# <cols> definitions which go out to column J) while only row 12 stays
# as the last populated row.
$ws.Range("J12").Font.Bold = $false

# --- Column width tweaks ---
$ws.Columns.Item(2).ColumnWidth = 13.166666666666666
$ws.Columns.Item(5).ColumnWidth = 1.1666666666666665
